$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# Update the "Date" metadata value
$wsMetadata.Range("B8").Value = "2026-01-16T13:49:34+00:00"

# Update the "Description" metadata value (drop the leading "Entrée ")
$wsMetadata.Range("B12").Value = "Quantité de produit"

# The "fr-lm-quantite-produit" base row's Definition cell shared the exact same
# string as the old Description value; keep it in sync so it still shares the
# (now-updated) string instead of leaving an orphaned duplicate behind.
$wsElements.Range("M2").Value = "Quantité de produit"
